$d = $word.ActiveDocument

# --- 1. Title (Heading1) ---
$d.Content.Find.Execute(
    "Review 157: [Short] End-to-End Speech Recognition Contextualization with Large Language Models, 30.09.23",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Review 156: Short : LONGLORA: EFFICIENT FINE-TUNING OF LONG CONTEXT LARGE LANGUAGE MODELS, 28.09.2023",
    2)

# --- 2. Paper link (bold) ---
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.10917v1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2309.12307v3",
    2)

# --- 3. Remove the now-superfluous empty paragraphs (work from the bottom up
#        so earlier indices stay valid) before rewriting the surviving ones. ---
# Before state (1-indexed): P1 Heading, P2 Paper, P3 blank, P4 hf-link,
# P5/P6/P7 blank, P8 audio-text, P9 blank, P10 how-they-did-it, P11 blank,
# P12 then-they-take, P13 blank, P14 blank.
$d.Paragraphs(12).Range.Delete()   # "לאחר מכן לוקחים..."
$d.Paragraphs(11).Range.Delete()   # blank
$d.Paragraphs(10).Range.Delete()   # "איך הם עשו זאת..."
$d.Paragraphs(9).Range.Delete()    # blank
$d.Paragraphs(7).Range.Delete()    # blank
$d.Paragraphs(6).Range.Delete()    # blank

# Remaining paragraphs now: P1 Heading, P2 Paper, P3 blank, P4 hf-link,
# P5 blank (was P5), P6 audio-text (was P8), P7 blank (was P13), P8 blank (was P14)

# --- 4. Paper preview link, now with four manual line breaks appended ---
$p4 = $d.Paragraphs(4)
$p4.Range.Text = "https://arxiv.org/abs/2309.12307.pdf"
$p4.Range.InsertAfter([char]11)
$p4.Range.InsertAfter([char]11)
$p4.Range.InsertAfter([char]11)
$p4.Range.InsertAfter([char]11)

# --- 5. First body paragraph: two Hebrew sentences joined by a blank line ---
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "כל מי שעוסק במודלי שפה בטח שמע על סוגיית אורך ההקשר (context length). אנו רוצים שהמודלים שלנו יהיו מסוגלים ״להחזיק בבטן״ כמות כמה שיותר גדולה של טקסט. אולם הקשר ארוך דורש כמות עצומה של משאבים לאימון ולאינפרנס. "
$p5.Range.InsertAfter([char]11)
$p5.Range.InsertAfter([char]11)
$p5.Range.InsertAfter("אז היום ב-#shorthebrewpapereviews אנו סוקרים מאמר שמציע גישת טיוב(fine-tuning) שמגדילה את אורך הקשר של מודל שפה. כלומר אם מודל שפה היה מאומן באימון מקדים (pretraining) עם אורך הקשר של 2048, השיטה המוצעת מאפשרת להאריכו פי 4 ל-8192. כמו שאתם רואים השם של השיטה מכיל את המילה LoRA שהיא שיטה מאוד פופולרית לפיין-טיון של מודלי שפה. ")

# --- 6. Second body paragraph: leading break, two sentences, blank line, closing break ---
$p6 = $d.Paragraphs(6)
$p6.Range.Text = [char]11
$p6.Range.InsertAfter("במקום לכייל (לשנות) את כל המשקלים של מודל השפה המכויל LoRA מעדכנת רק את התוספת למשקלי המודל (כמו ResNet). בנוסף התוספת למשקלים מיוצגת על ידי מטריצה עם רנק נמוך שניתן לתאר אותה על ידי מכפלה של מטריצות בעלות מימד נמוך יחסית. ")
$p6.Range.InsertAfter([char]11)
$p6.Range.InsertAfter([char]11)
$p6.Range.InsertAfter("אז מה מציע LongLoRA בנוסף? כדי להגדיל את אורך הקשר נגיד מ- 2048 ל-8192 היא מחלקת את 8192 טוקנים ל-4 קבוצות בעלות 2048 טוקנים כל אחת שעבור כל אחת מהם ציוני ה-attention מחושבים בנפרד (חיסכון פי 16 בחישובים). את זה עושים בחצי מהראשים. בשאר הראשים פשוט מזיזים את הקבוצות האלו בחצי גודל כלומר הקבוצה הראשונה תכילי טוקנים מ-1024 עד 3071, השניה מ 3072 ל 5195 וכדומה. טריק פשוט מאוד אבל מביא תוצאות לא רעות בכלל.")
$p6.Range.InsertAfter([char]11)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
